$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q1").Value = "wtkappa.scale_trim_round"
$ws.Range("E2").Value = 0.02324751973535285
$ws.Range("F2").Value = -0.01054416404922412
$ws.Range("Q2").Value = 0.782122905027933
